$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# @@ -727,25 +727,25 @@
$ws.Range("H2").Value = 92.2
$ws.Range("I2").Value = 89.75
$ws.Range("J2").Value = 102
$ws.Range("K2").Value = 89.75
$ws.Range("L2").Value = 102
$ws.Range("M2").Value = 23.25
$ws.Range("N2").Value = -328

# @@ -2062,25 +2062,22 @@
$ws.Range("H29").Value = 851.5
$ws.Range("I29").Value = 851.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2554.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2273.5
$ws.Range("N29").Value = $null

# @@ -2518,25 +2515,25 @@
$ws.Range("H38").Value = 2043.8518
$ws.Range("I38").Value = 348.7857
$ws.Range("J38").Value = 3869.3076
$ws.Range("K38").Value = 1046.3571
$ws.Range("L38").Value = 11607.9228
$ws.Range("M38").Value = -674.3571000000002
$ws.Range("N38").Value = -12351.9228

# @@ -3534,25 +3531,25 @@
$ws.Range("H58").Value = 897.3043
$ws.Range("I58").Value = 385.7143
$ws.Range("J58").Value = 1121.125
$ws.Range("K58").Value = 1157.1429
$ws.Range("L58").Value = 3363.375
$ws.Range("M58").Value = -1007.1429
$ws.Range("N58").Value = -3663.375

# @@ -4339,25 +4336,25 @@
$ws.Range("H74").Value = 3525
$ws.Range("I74").Value = 3466.6667
$ws.Range("J74").Value = 3700
$ws.Range("K74").Value = 3466.6667
$ws.Range("L74").Value = 3700
$ws.Range("M74").Value = -2530.6667
$ws.Range("N74").Value = -5572

# @@ -4437,22 +4434,22 @@
$ws.Range("H76").Value = 9262542
$ws.Range("I76").Value = 18520518
$ws.Range("K76").Value = 18520518
$ws.Range("M76").Value = -18520203

# @@ -4489,25 +4486,25 @@
$ws.Range("H77").Value = 3525
$ws.Range("I77").Value = 3466.6667
$ws.Range("J77").Value = 3700
$ws.Range("K77").Value = 17333.3335
$ws.Range("L77").Value = 18500
$ws.Range("M77").Value = -12653.3335
$ws.Range("N77").Value = -27860

# @@ -4587,22 +4584,22 @@
$ws.Range("H79").Value = 9262542
$ws.Range("I79").Value = 18520518
$ws.Range("K79").Value = 18520518
$ws.Range("M79").Value = -18519426

# @@ -4997,25 +4994,25 @@
$ws.Range("H87").Value = 11103.205
$ws.Range("J87").Value = 11421.7295
$ws.Range("L87").Value = 11421.7295
$ws.Range("N87").Value = -13917.7295

# @@ -5153,25 +5150,25 @@
$ws.Range("H90").Value = 11103.205
$ws.Range("J90").Value = 11421.7295
$ws.Range("L90").Value = 34265.1885
$ws.Range("N90").Value = -46745.1885

# @@ -7115,25 +7112,25 @@
$ws.Range("H129").Value = 1034.7142
$ws.Range("J129").Value = 1299.6
$ws.Range("L129").Value = 3898.8
$ws.Range("N129").Value = -13898.8

# @@ -7268,25 +7265,25 @@
$ws.Range("H132").Value = 2346.52
$ws.Range("I132").Value = 2142.1
$ws.Range("J132").Value = 3164.2
$ws.Range("K132").Value = 6426.299999999999
$ws.Range("L132").Value = 9492.599999999999
$ws.Range("M132").Value = -3896.299999999999
$ws.Range("N132").Value = -14552.6

# @@ -7421,25 +7418,25 @@
$ws.Range("H135").Value = 4404
$ws.Range("I135").Value = 3275.6365
$ws.Range("J135").Value = 12678.667
$ws.Range("K135").Value = 29480.7285
$ws.Range("L135").Value = 114108.003
$ws.Range("M135").Value = -26945.7285
$ws.Range("N135").Value = -119178.003

# @@ -7574,25 +7571,25 @@
$ws.Range("H138").Value = 4627.08
$ws.Range("I138").Value = 2405.4443
$ws.Range("J138").Value = 5448.781
$ws.Range("K138").Value = 7216.3329
$ws.Range("L138").Value = 16346.343
$ws.Range("M138").Value = -2076.3329
$ws.Range("N138").Value = -26626.343

$ws = $wb.Worksheets.Item("ARM")
# @@ -9534,22 +9531,22 @@
$ws.Range("H37").Value = 8219
$ws.Range("J37").Value = 8219
$ws.Range("L37").Value = 8219
$ws.Range("N37").Value = -8765

$ws = $wb.Worksheets.Item("CRP")
# @@ -23013,25 +23010,25 @@
$ws.Range("H31").Value = 9778.324000000001
$ws.Range("I31").Value = 3130.439
$ws.Range("J31").Value = 18037.818
$ws.Range("K31").Value = 3130.439
$ws.Range("L31").Value = 18037.818
$ws.Range("M31").Value = -2835.439
$ws.Range("N31").Value = -18627.818

# @@ -23166,25 +23163,25 @@
$ws.Range("H34").Value = 9778.324000000001
$ws.Range("I34").Value = 3130.439
$ws.Range("J34").Value = 18037.818
$ws.Range("K34").Value = 3130.439
$ws.Range("L34").Value = 18037.818
$ws.Range("M34").Value = -2928.439
$ws.Range("N34").Value = -18441.818

# @@ -23977,22 +23974,22 @@
$ws.Range("H50").Value = 8635.200000000001
$ws.Range("J50").Value = 8635.200000000001
$ws.Range("L50").Value = 8635.200000000001
$ws.Range("N50").Value = -9885.200000000001

# @@ -24026,25 +24023,25 @@
$ws.Range("H51").Value = 8519
$ws.Range("J51").Value = 9398.75
$ws.Range("L51").Value = 9398.75
$ws.Range("N51").Value = -10870.75

# @@ -24176,22 +24173,19 @@
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null

# @@ -24427,22 +24421,22 @@
$ws.Range("H59").Value = 11277
$ws.Range("J59").Value = 11277
$ws.Range("L59").Value = 11277
$ws.Range("N59").Value = -13567

# @@ -24476,25 +24470,25 @@
$ws.Range("H60").Value = 6875.25
$ws.Range("J60").Value = 8252.25
$ws.Range("L60").Value = 8252.25
$ws.Range("N60").Value = -9274.25

# @@ -24528,25 +24522,25 @@
$ws.Range("H61").Value = 8519
$ws.Range("J61").Value = 9398.75
$ws.Range("L61").Value = 9398.75
$ws.Range("N61").Value = -10094.75

# @@ -24880,22 +24874,22 @@
$ws.Range("H68").Value = 17295
$ws.Range("J68").Value = 17295
$ws.Range("L68").Value = 17295
$ws.Range("N68").Value = -18793

# @@ -25030,22 +25024,22 @@
$ws.Range("H71").Value = 17295
$ws.Range("J71").Value = 17295
$ws.Range("L71").Value = 51885
$ws.Range("N71").Value = -59373

# @@ -25180,22 +25174,22 @@
$ws.Range("H74").Value = 17167.666
$ws.Range("J74").Value = 17167.666
$ws.Range("L74").Value = 17167.666
$ws.Range("N74").Value = -18915.666

# @@ -25324,22 +25318,22 @@
$ws.Range("H77").Value = 17167.666
$ws.Range("J77").Value = 17167.666
$ws.Range("L77").Value = 51502.99800000001
$ws.Range("N77").Value = -60238.99800000001

$ws = $wb.Worksheets.Item("CUL")
# @@ -34732,25 +34726,25 @@
$ws.Range("H122").Value = 519.6923
$ws.Range("I122").Value = 422.82144
$ws.Range("J122").Value = 766.2727
$ws.Range("K122").Value = 3805.39296
$ws.Range("L122").Value = 6896.454299999999
$ws.Range("M122").Value = -1355.39296
$ws.Range("N122").Value = -11796.4543

# @@ -35500,25 +35494,25 @@
$ws.Range("H137").Value = 4600.625
$ws.Range("I137").Value = 2027.5
$ws.Range("J137").Value = 5458.3335
$ws.Range("K137").Value = 6082.5
$ws.Range("L137").Value = 16375.0005
$ws.Range("M137").Value = -982.5
$ws.Range("N137").Value = -26575.0005

$ws = $wb.Worksheets.Item("GSM")
# @@ -39183,22 +39177,22 @@
$ws.Range("H70").Value = 31259606
$ws.Range("I70").Value = 51144720
$ws.Range("K70").Value = 51144720
$ws.Range("M70").Value = -51144450

# @@ -39327,22 +39321,22 @@
$ws.Range("H73").Value = 31259606
$ws.Range("I73").Value = 51144720
$ws.Range("K73").Value = 51144720
$ws.Range("M73").Value = -51143784

$ws = $wb.Worksheets.Item("WVR")
# @@ -52636,25 +52630,25 @@
$ws.Range("H62").Value = 5853.6
$ws.Range("I62").Value = 6222.222
$ws.Range("J62").Value = 5300.6665
$ws.Range("K62").Value = 6222.222
$ws.Range("L62").Value = 5300.6665
$ws.Range("M62").Value = -5598.222
$ws.Range("N62").Value = -6548.6665

# @@ -52789,25 +52783,25 @@
$ws.Range("H65").Value = 5853.6
$ws.Range("I65").Value = 6222.222
$ws.Range("J65").Value = 5300.6665
$ws.Range("K65").Value = 31111.11
$ws.Range("L65").Value = 26503.3325
$ws.Range("M65").Value = -27991.11
$ws.Range("N65").Value = -32743.3325

# @@ -53037,19 +53031,25 @@
$ws.Range("H70").Value = 10500
$ws.Range("I70").Value = 8000
$ws.Range("J70").Value = 13000
$ws.Range("K70").Value = 8000
$ws.Range("L70").Value = 13000
$ws.Range("M70").Value = -7685
$ws.Range("N70").Value = -13630

# @@ -53178,19 +53178,25 @@
$ws.Range("H73").Value = 10500
$ws.Range("I73").Value = 8000
$ws.Range("J73").Value = 13000
$ws.Range("K73").Value = 8000
$ws.Range("L73").Value = 13000
$ws.Range("M73").Value = -6908
$ws.Range("N73").Value = -15184
